$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume number and week-covering dates)
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# Crime-statistics table updates (rows 14-31)
$ws.Range("M14").Value = -77.777777777777
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = 38.461538461538
$ws.Range("L15").Value = 63.636363636363
$ws.Range("M15").Value = 5.882352941176
$ws.Range("N15").Value = -41.935483870967
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 20
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = 6.493506493506
$ws.Range("L16").Value = -14.583333333333
$ws.Range("M16").Value = -47.096774193548
$ws.Range("N16").Value = -85.486725663716
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -8.333333333333
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 17.142857142857
$ws.Range("I17").Value = 255
$ws.Range("J17").Value = 281
$ws.Range("K17").Value = -9.252669039145
$ws.Range("L17").Value = -0.778210116731
$ws.Range("M17").Value = 40.109890109890
$ws.Range("N17").Value = -42.567567567567
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = -8.196721311475
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -64.331210191082
$ws.Range("N18").Value = -94.620557156580
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 233.333333333333
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 184.615384615385
$ws.Range("I19").Value = 186
$ws.Range("J19").Value = 229
$ws.Range("K19").Value = -18.777292576419
$ws.Range("L19").Value = -18.777292576419
$ws.Range("M19").Value = -19.480519480519
$ws.Range("N19").Value = -42.236024844720
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -9.090909090909
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = -55.263157894736
$ws.Range("L20").Value = -53.424657534246
$ws.Range("M20").Value = -70.689655172413
$ws.Range("N20").Value = -95.572916666666
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 38.095238095238
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = 34.567901234567
$ws.Range("I21").Value = 633
$ws.Range("J21").Value = 738
$ws.Range("K21").Value = -14.227642276422
$ws.Range("L21").Value = -16.600790513834
$ws.Range("M21").Value = -26.989619377162
$ws.Range("N21").Value = -80.106851037083
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 125
$ws.Range("I23").Value = 40
$ws.Range("K23").Value = -11.111111111111
$ws.Range("L23").Value = -28.571428571428
$ws.Range("M23").Value = 53.846153846153
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -5.555555555555
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = -10
$ws.Range("I24").Value = 717
$ws.Range("J24").Value = 734
$ws.Range("K24").Value = -2.316076294277
$ws.Range("L24").Value = 3.017241379310
$ws.Range("M24").Value = -11.699507389162
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 140
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 19.230769230769
$ws.Range("I25").Value = 375
$ws.Range("J25").Value = 324
$ws.Range("K25").Value = 15.740740740740
$ws.Range("L25").Value = 38.376383763837
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -20
$ws.Range("G26").Value = 67
$ws.Range("H26").Value = 4.477611940298
$ws.Range("I26").Value = 436
$ws.Range("J26").Value = 435
$ws.Range("K26").Value = 0.229885057471
$ws.Range("L26").Value = 7.125307125307
$ws.Range("M26").Value = -36.627906976744
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 25
$ws.Range("K27").Value = 8.695652173913
$ws.Range("L27").Value = 56.25
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -9.090909090909
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = 3.921568627450
$ws.Range("L28").Value = 3.921568627450
$ws.Range("C29").Value = "0"
$ws.Range("M29").Value = -78.947368421052
$ws.Range("N29").Value = -93.442622950819
$ws.Range("C30").Value = "0"
$ws.Range("M30").Value = -76.470588235294
$ws.Range("N30").Value = -92.156862745098
$ws.Range("F31").Value = 1
